# Edit script: adds camera data rows + a new summary column, and updates the
# sheet view (freeze header row) to match the authored commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

# ---------------------------------------------------------------------------
# 1. Fill in missing Locality / Date cells for the existing rows 9-16
# ---------------------------------------------------------------------------
$ws.Range("A10:A16").Value = "CojoHQ"
$ws.Range("C9:C16").Value = 45727
$ws.Range("C9:C16").NumberFormat = "d-mmm-yy"

# ---------------------------------------------------------------------------
# 2. New header cell for column S
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# 3. New data rows 17-26
#    (string-valued cells are written in the same order they were first
#    authored so the generated shared-string table lines up with the
#    original edit)
# ---------------------------------------------------------------------------

# --- Row 17 ---
$ws.Range("A17").Value = "PC "
$ws.Range("B17").Value = 14
$ws.Range("C17").Value = 45731
$ws.Range("D17").Value = "before 8:14"
$ws.Range("G17").Value = 0.40555555555555556
$ws.Range("H17").Value = 6303
$ws.Range("J17").Value = 0.55138888888888893
$ws.Range("K17").Value = 6324
$ws.Range("F19").Value = "new flower"
$ws.Range("L17").Value = "1st petal abscises"
$ws.Range("M17").NumberFormat = "h:mm"

$ws.Range("S1").Value = "Time available"

# --- Row 18 ---
$ws.Range("A18").Value = "PC "
$ws.Range("B18").Value = 15
$ws.Range("C18").Value = 45731
$ws.Range("D18").Value = 0.45416666666666666
$ws.Range("E18").Value = 6310
$ws.Range("G18").Value = 0.53055555555555556
$ws.Range("H18").Value = 6321
$ws.Range("J18").Value = 0.67638888888888893
$ws.Range("K18").Value = 6342
$ws.Range("M18").Value = 0.77361111111111114
$ws.Range("N18").Value = 6356

# --- Row 19 ---
$ws.Range("A19").Value = "PC "
$ws.Range("B19").Value = 16
$ws.Range("C19").Value = 45731
$ws.Range("D19").Value = 0.46805555555555556
$ws.Range("G19").Value = 0.53055555555555556
$ws.Range("H19").Value = 6321
$ws.Range("J19").Value = 0.63472222222222219
$ws.Range("K19").Value = 6336
$ws.Range("M19").Value = 0.73888888888888893
$ws.Range("N19").Value = 6351

# --- Row 20 ---
$ws.Range("A20").Value = "PC "
$ws.Range("B20").Value = 15
$ws.Range("C20").Value = 45732
$ws.Range("D20").Value = "NA"
$ws.Range("G20").Value = "NA"
$ws.Range("J20").Value = 0.69027777777777777
$ws.Range("K20").Value = 6434
$ws.Range("M20").Value = 0.73888888888888893
$ws.Range("N20").Value = 6441

# --- Row 21 ---
$ws.Range("A21").Value = "PC "
$ws.Range("B21").Value = 16
$ws.Range("C21").Value = 45732
$ws.Range("D21").Value = "NA"
$ws.Range("G21").Value = "NA"
$ws.Range("J21").Value = 0.64861111111111114
$ws.Range("K21").Value = 6428
$ws.Range("M21").Value = 0.71111111111111114
$ws.Range("N21").Value = 6437

# --- Row 22 ---
$ws.Range("A22").Value = "PC "
$ws.Range("B22").Value = 17
$ws.Range("C22").Value = 45732
$ws.Range("D22").Value = 0.35416666666666669
$ws.Range("E22").Value = 4133

# --- Row 23 ---
$ws.Range("A23").Value = "PC "
$ws.Range("B23").Value = 18
$ws.Range("C23").Value = 45731
$ws.Range("D23").Value = 0.3611111111111111
$ws.Range("E23").Value = 3991
$ws.Range("G23").Value = 0.4236111111111111
$ws.Range("H23").Value = 4000
$ws.Range("J23").Value = 0.63888888888888884
$ws.Range("K23").Value = 4031
$ws.Range("L23").Value = "1st petal abcises"

# --- Row 24 ---
$ws.Range("A24").Value = "PC "
$ws.Range("B24").Value = 19
$ws.Range("C24").Value = 45731
$ws.Range("D24").Value = 0.375
$ws.Range("E24").Value = 3993
$ws.Range("G24").Value = 0.47916666666666669
$ws.Range("H24").Value = 4008
$ws.Range("J24").Value = 0.66666666666666663
$ws.Range("K24").Value = 4035
$ws.Range("L24").Value = 0.70833333333333337
$ws.Range("M24").Value = 0.70833333333333337
$ws.Range("N24").Value = 4041
$ws.Range("N24").NumberFormat = "h:mm"

# --- Row 25 ---
$ws.Range("A25").Value = "PC "
$ws.Range("B25").Value = 20
$ws.Range("C25").Value = 45731
$ws.Range("D25").Value = 0.41666666666666669
$ws.Range("E25").Value = 3999
$ws.Range("G25").Value = 0.54861111111111116
$ws.Range("H25").Value = 4016
$ws.Range("J25").Value = "15;40"
$ws.Range("K25").Value = 4033
$ws.Range("L25").Value = 0.70138888888888884
$ws.Range("M25").Value = 0.70138888888888884

# --- Row 26 ---
$ws.Range("A26").Value = "PC "
$ws.Range("B26").Value = 21
$ws.Range("C26").Value = 45731
$ws.Range("D26").Value = 0.41666666666666669
$ws.Range("E26").Value = 3999
$ws.Range("F26").Value = "new flower"
$ws.Range("G26").Value = 0.54166666666666663
$ws.Range("H26").Value = 4017
$ws.Range("J26").Value = 0.63888888888888884
$ws.Range("K26").Value = 4031
$ws.Range("L26").Value = 0.71527777777777779
$ws.Range("M26").Value = 0.71527777777777779
$ws.Range("N26").Value = 4042

$ws.Range("N25").Value = 4040

# ---------------------------------------------------------------------------
# 4. Number formats for the new date (column C) and time (columns D/G/J/M/L)
#    cells in rows 17-26
# ---------------------------------------------------------------------------
$ws.Range("C17:C26").NumberFormat = "d-mmm-yy"

$ws.Range("D18:D19").NumberFormat = "h:mm"
$ws.Range("G17:G19").NumberFormat = "h:mm"
$ws.Range("J17:J21").NumberFormat = "h:mm"
$ws.Range("M18:M21").NumberFormat = "h:mm"

$ws.Range("D22:D26").NumberFormat = "h:mm"
$ws.Range("G23:G26").NumberFormat = "h:mm"
$ws.Range("J23:J24").NumberFormat = "h:mm"
$ws.Range("J26").NumberFormat = "h:mm"
$ws.Range("L24:L26").NumberFormat = "h:mm"
$ws.Range("M24:M26").NumberFormat = "h:mm"
$ws.Range("L25").NumberFormat = "h:mm"

# ---------------------------------------------------------------------------
# 5. Sheet view: freeze the header row and set the final selection
# ---------------------------------------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("P23").Select()
